$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change header from "from_number" to "roll_number"
$ws.Range("A1").Value = "roll_number"

# Update the active selection (as seen in sheetView) from J12 to A5
$ws.Range("A5").Select()
